$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harmoniza os nomes das colunas do cabecalho em caixa alta
$ws.Range("A1").Value = "FONTE_STN_COD"
$ws.Range("B1").Value = "FONTE_STN_DESCRICAO"
$ws.Range("C1").Value = "INTERPRETACAO"
$ws.Range("D1").Value = "DT_INICIO_VIGENCIA"
$ws.Range("E1").Value = "DT_FIM_VIGENCIA"

# Remove a formatacao numerica explicita (mas redundante) da coluna A
# de dados, restaurando o estilo padrao "Normal"
$ws.Range("A2:A86").Style = "Normal"
